$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").Value = 3.62
$ws.Range("E6").Value = 4.315
$ws.Range("F6").Value = 5.145
$ws.Range("I6").Value = 3.885
